$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new row 5 with "2021年" data, matching the formatting of row 4 (A2:A4 style)
$ws.Range("A4").Copy()
$ws.Range("A5").PasteSpecial(-4122)
$ws.Range("A5").Value = "2021年"

$ws.Range("B5").Value = 22.642
$ws.Range("C5").Value = 30.736
$ws.Range("D5").Value = 17.426
$ws.Range("E5").Value = 18.722
$ws.Range("F5").Value = 24.643
$ws.Range("G5").Value = 28.777
$ws.Range("H5").Value = 26.881
$ws.Range("I5").Value = 30.822
$ws.Range("J5").Value = 25.556
$ws.Range("K5").Value = 21.971
$ws.Range("L5").Value = 48.98
$ws.Range("M5").Value = 15.789
